$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the merge on B3:B5 (kept individual cell content/styles as-is)
# ---------------------------------------------------------------------------
$ws.Range("B3:B5").UnMerge()

# ---------------------------------------------------------------------------
# 2. Rows 6 & 7 no longer carry an explicit (taller) row height - AutoFit
#    restores the default height and drops the explicit ht attribute.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()

# ---------------------------------------------------------------------------
# 3. Rows 21-24 used a stray "applyFill" flavour of the bordered style (s=7);
#    normalise them back to the plain bordered style used everywhere else by
#    copying the format from A20 (style index 2) onto the affected cells.
# ---------------------------------------------------------------------------
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Fill the previously-empty row 25 and append the new observation rows
#    26-33 (content + matching bordered formatting copied from row 20).
# ---------------------------------------------------------------------------
$ws.Range("A20:C20").Copy()
$ws.Range("A25:C32").PasteSpecial(-4122)
$ws.Range("A33:C33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 33's middle cell (observation text) needs the word-wrap style like the
# other long multi-line observations (e.g. B6/B7) - copy that format too.
$ws.Range("B6").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A25").Value = "annual_inc_cat, loan_status"
$ws.Range("B25").Value = "More Charged off in 30k to 60k range"
$ws.Range("C25").Value = "Countplot"

$ws.Range("A26").Value = "purpose, loan_status,annual_inc"
$ws.Range("B26").Value = "Annual income with 50k to 60k with purpose home improvement have more charged off"
$ws.Range("C26").Value = "barplot"

$ws.Range("A27").Value = "int_rate_cat,annual_inc,loan_status"
$ws.Range("B27").Value = "Interest rate 20 to 25 percentage  with annual income 60000+ have more charged off"
$ws.Range("C27").Value = "barplot"

$ws.Range("A28").Value = "purpose, loan_status, int_rate"
$ws.Range("B28").Value = "Loan taken for house purpose has more charged off"
$ws.Range("C28").Value = "boxplot"

$ws.Range("A29").Value = "int_rate, grade, loan_status"
$ws.Range("B29").Value = "Interest high with grade G has more defaulters"
$ws.Range("C29").Value = "boxplot"

$ws.Range("A30").Value = "home_ownership, loan_status,annual_inc"
$ws.Range("B30").Value = "Mortgage home ownership with 60000 annual income have more charged off"
$ws.Range("C30").Value = "barplot"

$ws.Range("A31").Value = "annual_inc_cat, loan_status, int_rate"
$ws.Range("B31").Value = "Annual income with 90k-120k having high interest rate have more charged off"
$ws.Range("C31").Value = "barplot"

$ws.Range("A32").Value = "issue_year, loan_status"
$ws.Range("B32").Value = "In 2011 Charged off loans are more as per the above graph"
$ws.Range("C32").Value = "Countplot"

$ws.Range("A33").Value = "Correlation Matrix between columns"
$ws.Range("B33").Value = "positive correlation at revol_util and int_rate`npositive correlation at loan amount approved ratio and year issued"
$ws.Range("C33").Value = "Correlation Matrix"

# Undo the implicit row auto-height bump that typing wrapped, multi-line
# text into row 33 triggers, so it keeps the sheet's default row height
# (matching every other untouched row).
$ws.Rows.Item(33).AutoFit()

# ---------------------------------------------------------------------------
# 5. Scroll / selection bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B39").Select()
